$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 518; this shifts the former rows 518-536
# down to 519-537, matching the canonical diff.
$ws.Rows.Item(518).Insert()

# Populate the newly inserted row 518 with the new weekly price record.
$ws.Cells.Item(518, 1).Value = 4
$ws.Cells.Item(518, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(518, 3).Value = "Los Lagos"
$ws.Cells.Item(518, 4).Value = 45075
$ws.Cells.Item(518, 5).Value = 10
$ws.Cells.Item(518, 6).Value = 100112023
$ws.Cells.Item(518, 7).Value = "Brócoli"
$ws.Cells.Item(518, 8).Value = "Sin especificar"
$ws.Cells.Item(518, 9).Value = "Primera"
$ws.Cells.Item(518, 10).Value = 250
$ws.Cells.Item(518, 11).Value = 1700
$ws.Cells.Item(518, 12).Value = 1700
$ws.Cells.Item(518, 13).Value = 1700
$ws.Cells.Item(518, 14).Value = "$/unidad"
$ws.Cells.Item(518, 15).Value = "Región Metropolitana"
$ws.Cells.Item(518, 16).Value = 1700
$ws.Cells.Item(518, 17).Value = 1
$ws.Cells.Item(518, 18).Value = "Hortaliza"
